{"js": "// Update the worksheet date and the 25 division problems to the new\n// values, matching the commit's regenerated content.\nconst replacements = [\n  [\"2026-01-25 Sunday\", \"2026-01-26 Monday\"],\n  [\"645\u00f75=\", \"637\u00f75=\"],\n  [\"716\u00f75=\", \"735\u00f72=\"],\n  [\"682\u00f75=\", \"868\u00f72=\"],\n  [\"290\u00f78=\", \"577\u00f77=\"],\n  [\"946\u00f74=\", \"885\u00f73=\"],\n  [\"652\u00f72=\", \"638\u00f78=\"],\n  [\"664\u00f76=\", \"599\u00f79=\"],\n  [\"271\u00f79=\", \"796\u00f74=\"],\n  [\"619\u00f75=\", \"666\u00f73=\"],\n  [\"129\u00f77=\", \"566\u00f74=\"],\n  [\"269\u00f73=\", \"303\u00f72=\"],\n  [\"590\u00f75=\", \"879\u00f75=\"],\n  [\"283\u00f79=\", \"470\u00f77=\"],\n  [\"740\u00f73=\", \"159\u00f77=\"],\n  [\"257\u00f78=\", \"414\u00f77=\"],\n  [\"141\u00f75=\", \"144\u00f74=\"],\n  [\"958\u00f74=\", \"507\u00f77=\"],\n  [\"667\u00f73=\", \"248\u00f72=\"],\n  [\"732\u00f75=\", \"547\u00f77=\"],\n  [\"333\u00f76=\", \"370\u00f77=\"],\n  [\"551\u00f72=\", \"167\u00f78=\"],\n  [\"968\u00f73=\", \"516\u00f72=\"],\n  [\"228\u00f72=\", \"879\u00f76=\"],\n  [\"686\u00f78=\", \"877\u00f79=\"],\n  [\"258\u00f76=\", \"322\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division problems to the new\n# values, matching the commit's regenerated content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-25 Sunday\", \"2026-01-26 Monday\"),\n    @(\"645\u00f75=\", \"637\u00f75=\"),\n    @(\"716\u00f75=\", \"735\u00f72=\"),\n    @(\"682\u00f75=\", \"868\u00f72=\"),\n    @(\"290\u00f78=\", \"577\u00f77=\"),\n    @(\"946\u00f74=\", \"885\u00f73=\"),\n    @(\"652\u00f72=\", \"638\u00f78=\"),\n    @(\"664\u00f76=\", \"599\u00f79=\"),\n    @(\"271\u00f79=\", \"796\u00f74=\"),\n    @(\"619\u00f75=\", \"666\u00f73=\"),\n    @(\"129\u00f77=\", \"566\u00f74=\"),\n    @(\"269\u00f73=\", \"303\u00f72=\"),\n    @(\"590\u00f75=\", \"879\u00f75=\"),\n    @(\"283\u00f79=\", \"470\u00f77=\"),\n    @(\"740\u00f73=\", \"159\u00f77=\"),\n    @(\"257\u00f78=\", \"414\u00f77=\"),\n    @(\"141\u00f75=\", \"144\u00f74=\"),\n    @(\"958\u00f74=\", \"507\u00f77=\"),\n    @(\"667\u00f73=\", \"248\u00f72=\"),\n    @(\"732\u00f75=\", \"547\u00f77=\"),\n    @(\"333\u00f76=\", \"370\u00f77=\"),\n    @(\"551\u00f72=\", \"167\u00f78=\"),\n    @(\"968\u00f73=\", \"516\u00f72=\"),\n    @(\"228\u00f72=\", \"879\u00f76=\"),\n    @(\"686\u00f78=\", \"877\u00f79=\"),\n    @(\"258\u00f76=\", \"322\u00f79=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
